$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 720.7143
$ws.Range("I38").Value = 162
$ws.Range("J38").Value = 1031.1111
$ws.Range("K38").Value = 486
$ws.Range("L38").Value = 3093.3333
$ws.Range("M38").Value = -114
$ws.Range("N38").Value = -3837.3333

$ws.Range("H39").Value = 127
$ws.Range("I39").Value = 92.5
$ws.Range("J39").Value = 219
$ws.Range("K39").Value = 277.5
$ws.Range("L39").Value = 657
$ws.Range("M39").Value = 18.5
$ws.Range("N39").Value = -1249

$ws.Range("H98").Value = 709.7778
$ws.Range("I98").Value = 677.04346
$ws.Range("J98").Value = 898
$ws.Range("K98").Value = 677.04346
$ws.Range("L98").Value = 898
$ws.Range("M98").Value = 820.95654
$ws.Range("N98").Value = -3894

$ws.Range("H122").Value = 709.7778
$ws.Range("I122").Value = 677.04346
$ws.Range("J122").Value = 898
$ws.Range("K122").Value = 2031.13038
$ws.Range("L122").Value = 2694
$ws.Range("M122").Value = 418.8696199999999
$ws.Range("N122").Value = -7594

$ws.Range("H127").Value = 1061.8334
$ws.Range("I127").Value = 341.3
$ws.Range("J127").Value = 1962.5
$ws.Range("K127").Value = 1023.9
$ws.Range("L127").Value = 5887.5
$ws.Range("M127").Value = 3936.1
$ws.Range("N127").Value = -15807.5

$ws.Range("H132").Value = 982452.4
$ws.Range("I132").Value = 2377.878
$ws.Range("J132").Value = 5447236
$ws.Range("K132").Value = 7133.634
$ws.Range("L132").Value = 16341708
$ws.Range("M132").Value = -4603.634
$ws.Range("N132").Value = -16346768

$ws.Range("H137").Value = 1251090.2
$ws.Range("I137").Value = 1515926.9
$ws.Range("J137").Value = 2575.0715
$ws.Range("K137").Value = 4547780.699999999
$ws.Range("L137").Value = 7725.2145
$ws.Range("M137").Value = -4545230.699999999
$ws.Range("N137").Value = -12825.2145

$ws.Range("H141").Value = 2022.25
$ws.Range("I141").Value = 1389.3414
$ws.Range("J141").Value = 5729.2856
$ws.Range("K141").Value = 4168.0242
$ws.Range("L141").Value = 17187.8568
$ws.Range("M141").Value = 1011.9758
$ws.Range("N141").Value = -27547.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 18906696
$ws.Range("I61").Value = 20855214
$ws.Range("J61").Value = 200925.6
$ws.Range("K61").Value = 20855214
$ws.Range("L61").Value = 200925.6
$ws.Range("M61").Value = -20855002
$ws.Range("N61").Value = -201349.6

$ws.Range("H74").Value = 6707823.5
$ws.Range("I74").Value = 7961219.5
$ws.Range("K74").Value = 7961219.5
$ws.Range("M74").Value = -7960345.5

$ws.Range("H77").Value = 6707823.5
$ws.Range("I77").Value = 7961219.5
$ws.Range("K77").Value = 39806097.5
$ws.Range("M77").Value = -39801729.5

$ws.Range("H122").Value = 3269757.5
$ws.Range("I122").Value = 1814.7931
$ws.Range("K122").Value = 5444.379300000001
$ws.Range("M122").Value = -2994.379300000001

$ws.Range("H132").Value = 47202.465
$ws.Range("I132").Value = 28489
$ws.Range("J132").Value = 122056.336
$ws.Range("K132").Value = 85467
$ws.Range("L132").Value = 366169.008
$ws.Range("M132").Value = -82937
$ws.Range("N132").Value = -371229.008

$ws.Range("H136").Value = 18906696
$ws.Range("I136").Value = 20855214
$ws.Range("J136").Value = 200925.6
$ws.Range("K136").Value = 62565642
$ws.Range("L136").Value = 602776.8
$ws.Range("M136").Value = -62563092
$ws.Range("N136").Value = -607876.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12726.454
$ws.Range("I86").Value = 16447.125
$ws.Range("J86").Value = 2804.6667
$ws.Range("K86").Value = 16447.125
$ws.Range("L86").Value = 2804.6667
$ws.Range("M86").Value = -15324.125
$ws.Range("N86").Value = -5050.6667

$ws.Range("H89").Value = 12726.454
$ws.Range("I89").Value = 16447.125
$ws.Range("J89").Value = 2804.6667
$ws.Range("K89").Value = 82235.625
$ws.Range("L89").Value = 14023.3335
$ws.Range("M89").Value = -76619.625
$ws.Range("N89").Value = -25255.3335

$ws.Range("H105").Value = 17243520
$ws.Range("I105").Value = 27779742
$ws.Range("K105").Value = 27779742
$ws.Range("M105").Value = -27777995

$ws.Range("H134").Value = 1960.9
$ws.Range("I134").Value = 1268.2858
$ws.Range("J134").Value = 3577
$ws.Range("K134").Value = 3804.8574
$ws.Range("L134").Value = 10731
$ws.Range("M134").Value = -1269.8574
$ws.Range("N134").Value = -15801

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 35716644
$ws.Range("I58").Value = 50002504
$ws.Range("K58").Value = 50002504
$ws.Range("M58").Value = -50002301

$ws.Range("H132").Value = 21897.244
$ws.Range("I132").Value = 1457.1316
$ws.Range("J132").Value = 92508.55
$ws.Range("K132").Value = 4371.3948
$ws.Range("L132").Value = 277525.65
$ws.Range("M132").Value = -1841.3948
$ws.Range("N132").Value = -282585.65

$ws.Range("H134").Value = 20774.127
$ws.Range("I134").Value = 1168.6061
$ws.Range("J134").Value = 50182.41
$ws.Range("K134").Value = 3505.8183
$ws.Range("L134").Value = 150547.23
$ws.Range("M134").Value = -970.8182999999999
$ws.Range("N134").Value = -155617.23

$ws.Range("H136").Value = 35716644
$ws.Range("I136").Value = 50002504
$ws.Range("K136").Value = 150007512
$ws.Range("M136").Value = -150004962

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 29283.03
$ws.Range("I87").Value = 28653.5
$ws.Range("J87").Value = 29366.967
$ws.Range("K87").Value = 85960.5
$ws.Range("L87").Value = 88100.901
$ws.Range("M87").Value = -84712.5
$ws.Range("N87").Value = -90596.901

$ws.Range("H90").Value = 29283.03
$ws.Range("I90").Value = 28653.5
$ws.Range("J90").Value = 29366.967
$ws.Range("K90").Value = 257881.5
$ws.Range("L90").Value = 264302.703
$ws.Range("M90").Value = -251641.5
$ws.Range("N90").Value = -276782.703

$ws.Range("H136").Value = 2745
$ws.Range("I136").Value = 2141.6667
$ws.Range("J136").Value = 3107
$ws.Range("K136").Value = 6425.000100000001
$ws.Range("L136").Value = 9321
$ws.Range("M136").Value = -1325.000100000001
$ws.Range("N136").Value = -19521

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1278
$ws.Range("I122").Value = 830
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 2490
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -40
$ws.Range("N122").Value = -10750

$ws.Range("H126").Value = 1528.1818
$ws.Range("I126").Value = 1201.1111
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 3603.3333
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1133.3333
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 58365.94
$ws.Range("I132").Value = 34513.332
$ws.Range("J132").Value = 201481.6
$ws.Range("K132").Value = 103539.996
$ws.Range("L132").Value = 604444.8
$ws.Range("M132").Value = -101009.996
$ws.Range("N132").Value = -609504.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2437.8462
$ws.Range("I40").Value = 2437.8462
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2437.8462
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2301.8462
$ws.Range("N40").ClearContents()

$ws.Range("H122").Value = 2976.6943
$ws.Range("I122").Value = 2555.4211
$ws.Range("J122").Value = 3447.5293
$ws.Range("K122").Value = 7666.263300000001
$ws.Range("L122").Value = 10342.5879
$ws.Range("M122").Value = -5216.263300000001
$ws.Range("N122").Value = -15242.5879

$ws.Range("H132").Value = 73861.07000000001
$ws.Range("I132").Value = 2970
$ws.Range("J132").Value = 168382.5
$ws.Range("K132").Value = 8910
$ws.Range("L132").Value = 505147.5
$ws.Range("M132").Value = -6380
$ws.Range("N132").Value = -510207.5

$ws.Range("H136").Value = 98353.57000000001
$ws.Range("I136").Value = 63332.812
$ws.Range("J136").Value = 210420
$ws.Range("K136").Value = 189998.436
$ws.Range("L136").Value = 631260
$ws.Range("M136").Value = -187448.436
$ws.Range("N136").Value = -636360

$ws.Range("H139").Value = 43951.5
$ws.Range("J139").Value = 43951.5
$ws.Range("L139").Value = 43951.5
$ws.Range("N139").Value = -54231.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 4140
$ws.Range("I28").Value = 2850
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 2850
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = -2502
$ws.Range("N28").Value = -5696

$ws.Range("H122").Value = 1650.1522
$ws.Range("I122").Value = 1187.9642
$ws.Range("J122").Value = 2369.111
$ws.Range("K122").Value = 3563.8926
$ws.Range("L122").Value = 7107.333
$ws.Range("M122").Value = -1113.8926
$ws.Range("N122").Value = -12007.333

$ws.Range("H132").Value = 52425.22
$ws.Range("I132").Value = 32643.686
$ws.Range("J132").Value = 167817.5
$ws.Range("K132").Value = 97931.058
$ws.Range("L132").Value = 503452.5
$ws.Range("M132").Value = -95401.058
$ws.Range("N132").Value = -508512.5

$ws.Range("H136").Value = 30302.014
$ws.Range("I136").Value = 18357.75
$ws.Range("J136").Value = 81754.234
$ws.Range("K136").Value = 55073.25
$ws.Range("L136").Value = 245262.702
$ws.Range("M136").Value = -52523.25
$ws.Range("N136").Value = -250362.702
